$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.023.19'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +5.57%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.577.66'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +6.93%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '504.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.98'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.993'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.575'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.595.03'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.56'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.104'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.341'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.22%  '
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.013.48'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.984.56'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.62'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000139'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.587.70'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.78'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '339.90'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '59.80'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.420'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.683.59'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.77%  '
$ws.Range('E27').Value = '  +3.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.992'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0853'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.42'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '155.52'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.17'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.56%  '
$ws.Range('E34').Value = '  +2.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.71'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.97'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.12%  '
$ws.Range('E37').Value = '  +5.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.853'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +27.19%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.843'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.99%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.47'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.17%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.76'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.97%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '297.75'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.25%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '35.63'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0566'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.619'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0996'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.992'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.72'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +10.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.05%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0234'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.61%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.035.39'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.07%  '
